# Applies the "Added exercise answers and modified main workshop docs" edit:
#  1. Every paragraph that is a numbered list item (or the title, or the
#     already-existing blank "ListParagraph" spacer paragraphs) gets explicit
#     single line-spacing / 0pt-after spacing (w:spacing w:after="0" w:line="240"
#     w:lineRule="auto").
#  2. Seven new blank paragraphs are inserted as visual spacers between
#     question items (some carry the ListParagraph style, some are plain).
#  3. "... with the value y - x" is split into extra runs (adding a couple of
#     blank-looking leading spaces as a separate run).
#  4. The package name "nlstools" is corrected to "Hmisc", and the trailing
#     " package" run is split into " " + "package".

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$enDash = [char]0x2013

function Set-SingleSpacingNoAfter($paraIndex) {
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle -> w:line="240" w:lineRule="auto"
    $p.Range.ParagraphFormat.SpaceAfter = 0        # -> w:after="0"
}

# ---------------------------------------------------------------------
# Step 1: apply the spacing tweak to every paragraph that keeps its place
# (title, the 10 numbered list items, and the 2 pre-existing blank
# "ListParagraph" spacer paragraphs). Paragraph indices are stable here
# because nothing is inserted/removed yet.
# ---------------------------------------------------------------------
$spacingTargets = @(1, 2, 3, 4, 5, 6, 32, 33, 34, 35, 57, 58, 59)
foreach ($idx in $spacingTargets) {
    Set-SingleSpacingNoAfter $idx
}

# ---------------------------------------------------------------------
# Step 2: text edits.
# ---------------------------------------------------------------------

# 2a. "Create a new variable z with the value y - x" -> split the trailing
# run into " with the value " + "  " + "y - x". This paragraph has no
# w:proofErr markers, so it is safe to round-trip its WordOpenXML and patch
# the one run we care about, then write the whole paragraph back (whole-
# paragraph InsertXML is the only replacement granularity that is safe).
$p4 = $d.Paragraphs.Item(4)
$full4 = $p4.Range.WordOpenXML
$s4 = $full4.IndexOf("<w:p ")
if ($s4 -lt 0) { $s4 = $full4.IndexOf("<w:p>") }
$e4 = $full4.IndexOf("</w:p>") + 6
$pxml4 = $full4.Substring($s4, $e4 - $s4)

$oldRun4 = "<w:r><w:t xml:space=`"preserve`"> with the value y $enDash x</w:t></w:r>"
$newRun4 = "<w:r><w:t xml:space=`"preserve`"> with the value </w:t></w:r>" +
           "<w:r><w:t xml:space=`"preserve`">  </w:t></w:r>" +
           "<w:r><w:t>y $enDash x</w:t></w:r>"
$pxml4 = $pxml4.Replace($oldRun4, $newRun4)
$pxml4 = $pxml4.Replace("<w:p ", "<w:p xmlns:w='$wNs' ")
$p4.Range.InsertXML($pxml4)

# 2b. "Install the nlstools package" -> "Install the Hmisc package", with the
# trailing run additionally split into " " + "package". This paragraph DOES
# contain w:proofErr markers around the package-name run, and those get
# silently dropped if the paragraph is round-tripped through WordOpenXML
# (it merges runs across them), so instead:
#   - swap the misspelled word via Find/Replace, which only touches the
#     text of the single run it matches and leaves proofErr siblings alone;
#   - split the final " package" run by rebuilding just that paragraph's
#     body by hand (keeping the w:p opening tag's original attributes,
#     which WordOpenXML does preserve faithfully).
$null = $d.Content.Find.Execute("nlstools", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "Hmisc", 2)

$pPkgIndex = $d.Paragraphs.Count
$pPkg = $d.Paragraphs.Item($pPkgIndex)
$fullPkg = $pPkg.Range.WordOpenXML
$sPkg = $fullPkg.IndexOf("<w:p ")
if ($sPkg -lt 0) { $sPkg = $fullPkg.IndexOf("<w:p>") }
$tagEndPkg = $fullPkg.IndexOf(">", $sPkg)
$openTagPkg = $fullPkg.Substring($sPkg, $tagEndPkg - $sPkg + 1)
$openTagPkg = $openTagPkg.Replace("<w:p ", "<w:p xmlns:w='$wNs' ")

$pkgXml = $openTagPkg +
    "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr>" +
    "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Install the </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:t>Hmisc</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:r><w:t>package</w:t></w:r>" +
    "</w:p>"
$pPkg.Range.InsertXML($pkgXml)

# ---------------------------------------------------------------------
# Step 3: insert the seven new blank spacer paragraphs. Working from the
# bottom of the document upward keeps all the not-yet-processed paragraph
# indices valid.
# ---------------------------------------------------------------------
function Insert-BlankAfter($paraIndex, [bool]$listStyle) {
    $p = $d.Paragraphs.Item($paraIndex)
    $insertPoint = $d.Range($p.Range.End, $p.Range.End)
    if ($listStyle) {
        $xml = "<w:p xmlns:w='$wNs'><w:pPr><w:pStyle w:val='ListParagraph'/>" +
               "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/></w:pPr></w:p>"
    } else {
        $xml = "<w:p xmlns:w='$wNs'><w:pPr>" +
               "<w:spacing w:after='0' w:line='240' w:lineRule='auto'/></w:pPr></w:p>"
    }
    $insertPoint.InsertXML($xml)
}

# (paragraph index, gets ListParagraph style?) - in descending order.
Insert-BlankAfter 58 $true    # after "Return the second column ..."
Insert-BlankAfter 34 $false   # after "Return the second row ..."
Insert-BlankAfter 33 $true    # after "Determine the dimensions ..."
Insert-BlankAfter 5  $false   # after "Return the value of z ..."
Insert-BlankAfter 4  $false   # after "Create a new variable z ..."
Insert-BlankAfter 3  $false   # after "Assign the value of 345 to y"
Insert-BlankAfter 2  $true    # after "Assign the value of 120 to x"

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
